# PROS-6581 - CCRU - new KPI tables and POS 2019
#
# The "KPI_Source" lookup table (repeated, one row set per POS sheet) gets
# four label/file-name refreshes for 2018 -> 2019 tooling, and the
# previously-active tab moves from "Pos 2018 - MT - Convenience Small"
# back to "Pos 2018 - FT".

$wb = $excel.ActiveWorkbook

# Update the shared lookup table cells (B7, B8, A11, C11) on every POS sheet.
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $ws.Range("B7").Value = "Contract Execution 2018"
    $ws.Range("B8").Value = "Equipment Execution 2018"
    $ws.Range("A11").Value = "KPI_CONVERSION"
    $ws.Range("C11").Value = "KPIs_2018/KPIConversion2018.xlsx"
}

# Restore the previously-selected cell on the sheet that used to be active
# ("Pos 2018 - MT - Convenience Small") before moving the active tab.
$wsConvenienceSmall = $wb.Worksheets.Item("Pos 2018 - MT - Convenience Small")
$wsConvenienceSmall.Activate()
$wsConvenienceSmall.Range("C11").Select() | Out-Null

# Make "Pos 2018 - FT" the active tab/selected cell (was tab index 2 / 0-based,
# now tab index 0 / 0-based).
$wsFT = $wb.Worksheets.Item("Pos 2018 - FT")
$wsFT.Activate()
$wsFT.Range("B17").Select() | Out-Null
